# Insert a new data row before row 23 (shifts existing rows 23-80 down to 24-81)
# and populate it with the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(23).Insert()

$ws.Cells.Item(23, 1).Value = 1
$ws.Cells.Item(23, 2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(23, 3).Value = 'Arica y Parinacota'
$ws.Cells.Item(23, 4).Value = 44791
$ws.Cells.Item(23, 5).Value = 15
$ws.Cells.Item(23, 6).Value = 100114001
$ws.Cells.Item(23, 7).Value = 'Papa'
$ws.Cells.Item(23, 8).Value = 'Asterix'
$ws.Cells.Item(23, 9).Value = '1a (cosecha lavada)'
$ws.Cells.Item(23, 10).Value = 1000
$ws.Cells.Item(23, 11).Value = 13000
$ws.Cells.Item(23, 12).Value = 14000
$ws.Cells.Item(23, 13).Value = 13500
$ws.Cells.Item(23, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(23, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(23, 16).Value = 540
$ws.Cells.Item(23, 17).Value = 25
$ws.Cells.Item(23, 18).Value = 'Hortaliza'
